# ProjectoSena.xlsx - "Add files via upload" edit
# Updates the OBJETIVOS sheet: retitles the general/specific objectives,
# trims the long specific-objectives list down to 3 items, and removes the
# now-unused leading blank row / trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OBJETIVOS")

# Drop the leading blank row (row 1 carried no data, only a height) so the
# sheet's used range starts at row 2 again, matching the uploaded file.
$ws.Rows.Item(1).AutoFit()

# --- Text updates -----------------------------------------------------
$ws.Range("B6").Value = "OBJETIVO GENERAL"
$ws.Range("C6").Value = "OBJETIVOS ESPECIFICOS"

$ws.Range("B7").Value = "Desarrollar un aplicativo para el agendamiento y seguimiento de servicios de mantenimiento electrico y electronico para la empresa meca ing"
$ws.Range("C7").Value = "Implementar un sistema de registro e inicio de sesión para usuarios con datos personales y credenciales (correo/contraseña)."

$ws.Range("C8").Value = "Notificar a técnicos y usuarios con recordatorios 24 horas antes de la cita."
$ws.Range("C9").Value = "Sincronizar las agendas entre usuarios y técnicos para evitar traslapes o duplicidad en los servicios."
$ws.Range("C10").Value = "Visualizar las citas asignadas a los técnicos para su conocimiento y ejecución."

$ws.Range("C9").HorizontalAlignment = -4130
$ws.Range("C10").HorizontalAlignment = -4130

# --- Row-height touch-ups to match the trimmed content -----------------
$ws.Rows.Item(8).RowHeight = 28
$ws.Rows.Item(9).RowHeight = 47

# Mirror the title row's bordered style a couple of columns further right
# (matches the formatting carried over from the original upload).
$ws.Range("B3").Copy()
$ws.Range("E2:F2").PasteSpecial(-4122)

# --- Remove the now-unused trailing objective rows ----------------------
# (this also shrinks the B7:B20 merge down to B7:B10 automatically)
$ws.Range("B11:C20").EntireRow.Delete()

# Leave the selection where the uploaded workbook had it.
$ws.Range("G10").Select()
